$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper cell used as a staging area so text that looks numeric (e.g. "1.017")
# round-trips through Copy + PasteSpecial(values) as literal text, instead of
# Excel auto-converting it to a number when assigned directly via .Value.
$helper = $ws.Range("Z1")

$updates = [ordered]@{
    'D2' = '27.494.10'
    'E2' = '  +2.19%  '
    'D3' = '1.872.38'
    'E3' = '  +1.60%  '
    'D4' = '1.017'
    'E4' = '  +0.86%  '
    'D5' = '312.79'
    'E5' = '  +1.07%  '
    'D6' = '1.015'
    'E6' = '  +0.79%  '
    'D7' = '0.4780'
    'E7' = '  +0.61%  '
    'E8' = '  +2.73%  '
    'D9' = '0.07372'
    'E9' = '  +2.39%  '
    'D10' = '0.9375'
    'E10' = '  +1.29%  '
    'D11' = '20.72'
    'E11' = '  +5.60%  '
    'D12' = '0.07831'
    'E12' = '  +2.35%  '
    'D13' = '1.867.77'
    'E13' = '  -0.04%  '
    'D14' = '5.443'
    'E14' = '  +2.51%  '
    'D15' = '6.584'
    'E15' = '  +2.87%  '
    'D16' = '90.88'
    'E16' = '  +2.49%  '
    'D17' = '1.018'
    'E17' = '  +0.88%  '
    'D18' = '0.000008901'
    'E18' = '  +3.12%  '
    'E19' = '  +0.77%  '
    'D20' = '14.91'
    'E20' = '  +2.53%  '
    'D21' = '27.537.89'
    'E21' = '  +2.26%  '
    'D22' = '5.130'
    'E22' = '  +1.60%  '
    'D23' = '10.72'
    'E23' = '  +0.72%  '
    'D24' = '1.962'
    'E24' = '  +2.06%  '
    'D25' = '154.12'
    'E25' = '  +1.22%  '
    'D26' = '18.54'
    'E26' = '  +2.20%  '
    'D27' = '2.017'
    'E27' = '  +0.72%  '
    'D28' = '115.92'
    'E28' = '  +1.41%  '
    'D29' = '4.992'
    'E29' = '  +0.94%  '
    'D30' = '0.08936'
    'E30' = '  +0.92%  '
    'D31' = '3.351'
    'E31' = '  +1.89%  '
    'D32' = '1.219'
    'E32' = '  +4.08%  '
    'E33' = '  +2.93%  '
    'D34' = '0.7516'
    'E34' = '  +0.44%  '
    'D35' = '2.689'
    'E35' = '  -2.36%  '
    'D36' = '0.02054'
    'E36' = '  +5.48%  '
    'D37' = '1.117'
    'E37' = '  +2.48%  '
    'D38' = '0.05303'
    'E38' = '  +0.82%  '
    'D39' = '3.007'
    'E39' = '  +1.55%  '
    'D40' = '0.5347'
    'E40' = '  +2.64%  '
    'D41' = '7.083'
    'E41' = '  +1.70%  '
    'D42' = '0.1527'
    'E42' = '  +0.91%  '
    'D43' = '8.420'
    'E43' = '  +2.51%  '
    'B44' = 'EnergySwap'
    'C44' = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
    'D44' = '10.66'
    'E44' = '  +1.69%  '
    'B45' = 'Decentraland'
    'C45' = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
    'D45' = '0.4824'
    'E45' = '  +2.09%  '
    'E46' = '  +0.76%  '
    'D47' = '1.662'
    'D48' = '102.95'
    'E48' = '  +1.18%  '
    'D49' = '67.27'
    'E49' = '  +2.59%  '
    'E50' = '  +1.09%  '
    'D51' = '0.9132'
    'E51' = '  +3.13%  '
}

foreach ($ref in $updates.Keys) {
    $helper.Formula = "=""" + $updates[$ref].Replace('"', '""') + """"
    $helper.Copy()
    $ws.Range($ref).PasteSpecial(-4163)
}

$helper.ClearContents()
$excel.CutCopyMode = 0
